$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.8712721145745577
$ws.Range("I3").Value = 0.0250905950965404
$ws.Range("K3").Value = 172.2083333333333

$ws.Range("Q3").Value = 16
$ws.Range("R3").Value = 25
$ws.Range("S3").Value = 95
$ws.Range("T3").Value = 182
$ws.Range("U3").Value = 227
$ws.Range("V3").Value = 5871
$ws.Range("W3").Value = 5862
$ws.Range("X3").Value = 5792
$ws.Range("Y3").Value = 5705
$ws.Range("Z3").Value = 5660

$ws.Range("AF3").Value = 0.997282
$ws.Range("AG3").Value = 0.995753
$ws.Range("AH3").Value = 0.983863
$ws.Range("AI3").Value = 0.9690839999999999
$ws.Range("AJ3").Value = 0.96144
